# Plano de Ação - update per commit: "Deletado (documento de mvv) e sensores,
# teve uma modificação no plano de ação."
# - Title row changed from "Plano de Ação do Projeto Lixeira Inteligente" to
#   "Plano de Ação do Projeto Smart trash"
# - Tasks reordered/renumbered, statuses mostly moved to "Concluído",
#   %concluido moved to 100%, deadlines (Prazo) filled in, "Gabreil" -> "Gabriel"
# - Removed the now-unused "Em andamento" status value (no row still uses it)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title (merged A1:G1)
$ws.Range("A1").Value = "Plano de Ação do Projeto Smart trash"

# Rows 4,5,7,8,9,10 never had a date in column D before, so their D cell is
# not yet formatted as a date. Grab the existing date format (already used by
# D3/D6) and spread it onto the rest of the D column before writing values,
# so Excel doesn't have to mint a brand-new number format for them.
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D4:D5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D7:D10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 3 - Tarefa 1
$ws.Range("A3").Value = "Tarefa 1 - Conferir Documentação"
$ws.Range("B3").Value = "Normal"
$ws.Range("C3").Value = "Concluído "
$ws.Range("D3").Value = (Get-Date -Year 2020 -Month 4 -Day 8).Date
$ws.Range("E3").Value = "Todos"
$ws.Range("F3").Value = 1

# Row 4 - Tarefa 2
$ws.Range("A4").Value = "Tarefa 2-Site Estático Dashboard (Google Charts)"
$ws.Range("B4").Value = "Alta"
$ws.Range("C4").Value = "Concluído "
$ws.Range("D4").Value = (Get-Date -Year 2020 -Month 5 -Day 1).Date
$ws.Range("E4").Value = "Todos"
$ws.Range("F4").Value = 1

# Row 5 - Tarefa 3
$ws.Range("A5").Value = "Tarefa 3-Site Estático Institucional - Local"
$ws.Range("B5").Value = "Alta"
$ws.Range("C5").Value = "Concluído "
$ws.Range("D5").Value = (Get-Date -Year 2020 -Month 5 -Day 1).Date
$ws.Range("E5").Value = "Rafael/Stefany"
$ws.Range("F5").Value = 1

# Row 6 - Tarefa 4
$ws.Range("A6").Value = "Tarefa 4-Planilha de BackLog / Planilha de Sprints"
$ws.Range("B6").Value = "Normal"
$ws.Range("C6").Value = "Concluído "
$ws.Range("D6").Value = (Get-Date -Year 2020 -Month 5 -Day 1).Date
$ws.Range("E6").Value = "Product owner"
$ws.Range("F6").Value = 1

# Row 7 - Tarefa 5
$ws.Range("A7").Value = "Tarefa 5-Especificação do Analytics"
$ws.Range("B7").Value = "Alta"
$ws.Range("C7").Value = "Não iniciada"
$ws.Range("D7").Value = (Get-Date -Year 2020 -Month 5 -Day 2).Date
$ws.Range("E7").Value = "Todos"
$ws.Range("F7").Value = 0

# Row 8 - Tarefa 6
$ws.Range("A8").Value = "Tarefa 6- Diagrama de Arquitetura Local (Arduíno)"
$ws.Range("B8").Value = "Alta"
$ws.Range("C8").Value = "Concluído "
$ws.Range("D8").Value = (Get-Date -Year 2020 -Month 5 -Day 3).Date
$ws.Range("E8").Value = "Yuri/Graziela"
$ws.Range("F8").Value = 1

# Row 9 - Tarefa 7
$ws.Range("A9").Value = "Tarefa 7- Tabelas criadas no Azure"
$ws.Range("B9").Value = "Alta"
$ws.Range("C9").Value = "Concluído "
$ws.Range("D9").Value = (Get-Date -Year 2020 -Month 5 -Day 4).Date
$ws.Range("E9").Value = "Bruno/Gabriel"
$ws.Range("F9").Value = 1

# Row 10 - Tarefa 8
$ws.Range("A10").Value = "Tarefa 8-Teste Integrado (Arduino+DB) + API local com Node.JS"
$ws.Range("B10").Value = "Alta"
$ws.Range("C10").Value = "Concluído "
$ws.Range("D10").Value = (Get-Date -Year 2020 -Month 5 -Day 5).Date
$ws.Range("E10").Value = "Todos"
$ws.Range("F10").Value = 1

# Update the selected cell shown when the sheet was last saved
$ws.Range("G4").Select()
